{"js": "// Change the semester reference from \"summer semester 2021\" to\n// \"winter semester 2021/2022\" (WiSe 21/22), as described by the commit\n// \"change text to WiSe21/22\".\n//\n// The surrounding sentence reads:\n//   \"...in the summer semester 2021. The students learned...\"\n// and becomes:\n//   \"...in the winter semester 2021/2022. The students learned...\"\n\nconst searchResults = context.document.body.search(\"summer semester 2021\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"summer semester 2021\" in the document body.');\n}\n\n// Replace the matched text in place, preserving the run's existing\n// character formatting (font, color, size) since insertText keeps the\n// formatting of the range it replaces.\nsearchResults.items[0].insertText(\"winter semester 2021/2022\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change the semester reference from \"summer semester 2021\" to\n# \"winter semester 2021/2022\" (WiSe 21/22), as described by the commit\n# \"change text to WiSe21/22\".\n#\n# The surrounding sentence reads:\n#   \"...in the summer semester 2021. The students learned...\"\n# and becomes:\n#   \"...in the winter semester 2021/2022. The students learned...\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"summer semester 2021\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"winter semester 2021/2022\"\n\n# wdFindContinue = 1, wdReplaceOne = 1\n$found = $find.Execute(\n    $find.Text,        # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    $find.Replacement.Text,  # ReplaceWith\n    1                  # Replace (wdReplaceOne)\n)\n\nif (-not $found) {\n    throw 'Could not find \"summer semester 2021\" in the document.'\n}\n"}
